# Enter police custody location
#
# The template has a legacy FORMTEXT field right after the question
# "If the offender is in police custody, state where:". That field has
# no default value configured, so Word renders its (empty) cached result
# as a run of five placeholder space characters. Replace that cached
# result with the {{custody_status_details}} merge-field placeholder.

$d = $word.ActiveDocument

# Find the question label that precedes the target form field.
$label = $d.Content
$found = $label.Find.Execute( `
    "If the offender is in police custody, state where:", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'police custody, state where' label"
}

# The form field's cached (empty) result is rendered as five consecutive
# "space-like" placeholder characters. Search for that block within the
# remainder of the document, starting right after the label.
$placeholderChar = [char]0x2002
$placeholderBlock = "$placeholderChar$placeholderChar$placeholderChar$placeholderChar$placeholderChar"

$fieldResult = $d.Range($label.End, $d.Content.End)
$foundBlock = $fieldResult.Find.Execute( `
    $placeholderBlock, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

if (-not $foundBlock) {
    throw "Could not find the cached form-field result to replace"
}

# Collapse the five placeholder runs down to the single merge-field run.
$fieldResult.Text = "{{custody_status_details}}"
